# ---------------------------------------------------------------------------
# webanalytics-pageload.xlsx maintenance edit:
#   - Add two new "blog" test-data sheets (BlogPostPage, BlogSeriesPage)
#   - BlogSeriesPage reuses/renames the previously-empty "Sheet2"
#   - BlogPostPage is a brand new sheet
#   - Re-order all sheets alphabetically
#   - Remove a stale/duplicate row from PDQPage
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Turn the empty "Sheet2" into "BlogSeriesPage" and fill it in -------
$blogSeries = $wb.Worksheets.Item("Sheet2")
$blogSeries.Name = "BlogSeriesPage"

$blogSeries.Range("A1").Value = "Path"
$blogSeries.Range("B1").Value = "ContentType"

$blogSeries.Range("A2").Value = "/news-events/cancer-currents-blog"
$blogSeries.Range("B2").Value = "Blog Series"

$blogSeries.Range("A3").Value = "/espanol/noticias/temas-y-relatos-blog"
$blogSeries.Range("B3").Value = "Blog Series"

$blogSeries.Range("A4").Value = "/research/key-initiatives/ras/ras-central/blog"
$blogSeries.Range("B4").Value = "Blog Series"

# copy the bold/shaded header formatting that every other sheet uses
$wb.Worksheets.Item("CTHPPage").Range("A1:B1").Copy()
$blogSeries.Range("A1:B1").PasteSpecial(-4122)

$blogSeries.Columns("A:B").AutoFit()
$blogSeries.Range("A5").Select()

# --- 2. Add a brand new "BlogPostPage" sheet and fill it in ----------------
$blogPost = $wb.Worksheets.Add($blogSeries)
$blogPost.Name = "BlogPostPage"

$blogPost.Range("A1").Value = "Path"
$blogPost.Range("B1").Value = "ContentType"

$blogPost.Range("A2").Value = "/news-events/cancer-currents-blog/2018/selumetinib-nf1-neurofibromas"
$blogPost.Range("B2").Value = "Blog Post"

$blogPost.Range("A3").Value = "/espanol/noticias/temas-y-relatos-blog/2018/selumetinib-neurofibromas-nf1"
$blogPost.Range("B3").Value = "Blog Post"

$blogPost.Range("A4").Value = "/about-nci/organization/cgh/blog/2017/cancer-research-day"
$blogPost.Range("B4").Value = "Blog Post"

$blogPost.Range("A5").Value = "/espanol/instituto/organizacion/salud-mundial/blog/2017/dia-mundial-investigacion"
$blogPost.Range("B5").Value = "Blog Post"

$wb.Worksheets.Item("CTHPPage").Range("A1:B1").Copy()
$blogPost.Range("A1:B1").PasteSpecial(-4122)

$blogPost.Columns("A:B").AutoFit()
$blogPost.Range("A6").Select()

# --- 3. Re-order the sheets: BlogPostPage, BlogSeriesPage, then the rest ---
# alphabetically, exactly as the workbook now lists them.
$blogPost.Move($wb.Worksheets.Item(1))
$blogSeries.Move($wb.Worksheets.Item(2))

$wb.Worksheets.Item("CTHPPage").Move($wb.Worksheets.Item(3))
$wb.Worksheets.Item("HomePage").Move($wb.Worksheets.Item(4))
$wb.Worksheets.Item("InnerPage").Move($wb.Worksheets.Item(5))
$wb.Worksheets.Item("LandingPage").Move($wb.Worksheets.Item(6))
$wb.Worksheets.Item("PDQPage").Move($wb.Worksheets.Item(7))
$wb.Worksheets.Item("TopicPage").Move($wb.Worksheets.Item(8))

# --- 4. Clean up PDQPage: drop the stale duplicate screening-overview row --
$pdq = $wb.Worksheets.Item("PDQPage")
$pdq.Rows("8:8").Delete()
$pdq.Range("A11").Select()

# --- 5. Leave PDQPage active/selected, matching the saved workbook view ----
$pdq.Select()
